$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 328-358: shorten/condense team names in column D ---
# (shared-string table gets regenerated/pruned by the engine once the old long-form
#  names are no longer referenced by any cell, matching the target sharedStrings.xml)
$ws.Range("D328").Value = "SanFrancisco"
$ws.Range("D332").Value = "NYGiants"
$ws.Range("D334").Value = "NewEngland"
$ws.Range("D342").Value = "TampaBay"
$ws.Range("D349").Value = "NYJets"
$ws.Range("D350").Value = "LasVegas"
$ws.Range("D352").Value = "LAChargers"
$ws.Range("D354").Value = "LARams"
$ws.Range("D355").Value = "KansasCity"
$ws.Range("D356").Value = "NewOrleans"
$ws.Range("D357").Value = "SanFrancisco"
$ws.Range("D358").Value = "GreenBay"

# --- Fill in new score/odds rows 360-389 (previously blank placeholder rows) ---
$ws.Rows(360).ClearFormats()
$ws.Range("A360").Value = 1128
$ws.Range("B360").Value = 275
$ws.Range("C360").Value = "V"
$ws.Range("D360").Value = "Pittsburgh"
$ws.Range("E360").Value = 3
$ws.Range("F360").Value = 13
$ws.Range("G360").Value = 0
$ws.Range("H360").Value = 8
$ws.Range("I360").Value = 24
$ws.Range("J360").Value = 40
$ws.Range("K360").Value = 39.5
$ws.Range("L360").Value = 115
$ws.Range("M360").Value = 19.5
$ws.Rows(361).ClearFormats()
$ws.Range("A361").Value = 1128
$ws.Range("B361").Value = 276
$ws.Range("C361").Value = "H"
$ws.Range("D361").Value = "Indianapolis"
$ws.Range("E361").Value = 0
$ws.Range("F361").Value = 3
$ws.Range("G361").Value = 14
$ws.Range("H361").Value = 0
$ws.Range("I361").Value = 17
$ws.Range("J361").Value = 3
$ws.Range("K361").Value = 2.5
$ws.Range("L361").Value = -135
$ws.Range("M361").Value = 3
$ws.Rows(362).ClearFormats()
$ws.Range("A362").Value = 1201
$ws.Range("B362").Value = 301
$ws.Range("C362").Value = "V"
$ws.Range("D362").Value = "Buffalo"
$ws.Range("E362").Value = 3
$ws.Range("F362").Value = 14
$ws.Range("G362").Value = 0
$ws.Range("H362").Value = 7
$ws.Range("I362").Value = 24
$ws.Range("J362").Value = 5.5
$ws.Range("K362").Value = 4
$ws.Range("L362").Value = -200
$ws.Range("M362").Value = 0.5
$ws.Rows(363).ClearFormats()
$ws.Range("A363").Value = 1201
$ws.Range("B363").Value = 302
$ws.Range("C363").Value = "H"
$ws.Range("D363").Value = "NewEngland"
$ws.Range("E363").Value = 7
$ws.Range("F363").Value = 0
$ws.Range("G363").Value = 0
$ws.Range("H363").Value = 3
$ws.Range("I363").Value = 10
$ws.Range("J363").Value = 45.5
$ws.Range("K363").Value = 44
$ws.Range("L363").Value = 175
$ws.Range("M363").Value = 21
$ws.Rows(364).ClearFormats()
$ws.Range("A364").Value = 1204
$ws.Range("B364").Value = 451
$ws.Range("C364").Value = "V"
$ws.Range("D364").Value = "NYJets"
$ws.Range("E364").Value = 3
$ws.Range("F364").Value = 3
$ws.Range("G364").Value = 6
$ws.Range("H364").Value = 10
$ws.Range("I364").Value = 22
$ws.Range("J364").Value = 42
$ws.Range("K364").Value = 43.5
$ws.Range("L364").Value = 130
$ws.Range("M364").Value = 20.5
$ws.Rows(365).ClearFormats()
$ws.Range("A365").Value = 1204
$ws.Range("B365").Value = 452
$ws.Range("C365").Value = "H"
$ws.Range("D365").Value = "Minnesota"
$ws.Range("E365").Value = 3
$ws.Range("F365").Value = 17
$ws.Range("G365").Value = 0
$ws.Range("H365").Value = 7
$ws.Range("I365").Value = 27
$ws.Range("J365").Value = 3
$ws.Range("K365").Value = 2.5
$ws.Range("L365").Value = -150
$ws.Range("M365").Value = 0.5
$ws.Rows(366).ClearFormats()
$ws.Range("A366").Value = 1204
$ws.Range("B366").Value = 453
$ws.Range("C366").Value = "V"
$ws.Range("D366").Value = "Denver"
$ws.Range("E366").Value = 3
$ws.Range("F366").Value = 3
$ws.Range("G366").Value = 3
$ws.Range("H366").Value = 0
$ws.Range("I366").Value = 9
$ws.Range("J366").Value = 41.5
$ws.Range("K366").Value = 40.5
$ws.Range("L366").Value = 330
$ws.Range("M366").Value = 19.5
$ws.Rows(367).ClearFormats()
$ws.Range("A367").Value = 1204
$ws.Range("B367").Value = 454
$ws.Range("C367").Value = "H"
$ws.Range("D367").Value = "Baltimore"
$ws.Range("E367").Value = 0
$ws.Range("F367").Value = 3
$ws.Range("G367").Value = 0
$ws.Range("H367").Value = 7
$ws.Range("I367").Value = 10
$ws.Range("J367").Value = 6
$ws.Range("K367").Value = 10
$ws.Range("L367").Value = -400
$ws.Range("M367").Value = 5
$ws.Rows(368).ClearFormats()
$ws.Range("A368").Value = 1204
$ws.Range("B368").Value = 455
$ws.Range("C368").Value = "V"
$ws.Range("D368").Value = "Pittsburgh"
$ws.Range("E368").Value = 3
$ws.Range("F368").Value = 13
$ws.Range("G368").Value = 3
$ws.Range("H368").Value = 0
$ws.Range("I368").Value = 19
$ws.Range("J368").Value = 42.5
$ws.Range("K368").Value = "pk"
$ws.Range("L368").Value = -125
$ws.Range("M368").Value = 21
$ws.Rows(369).ClearFormats()
$ws.Range("A369").Value = 1204
$ws.Range("B369").Value = 456
$ws.Range("C369").Value = "H"
$ws.Range("D369").Value = "Atlanta"
$ws.Range("E369").Value = 0
$ws.Range("F369").Value = 6
$ws.Range("G369").Value = 7
$ws.Range("H369").Value = 3
$ws.Range("I369").Value = 16
$ws.Range("J369").Value = 1
$ws.Range("K369").Value = 42.5
$ws.Range("L369").Value = 105
$ws.Range("M369").Value = 3
$ws.Rows(370).ClearFormats()
$ws.Range("A370").Value = 1204
$ws.Range("B370").Value = 457
$ws.Range("C370").Value = "V"
$ws.Range("D370").Value = "Tennessee"
$ws.Range("E370").Value = 7
$ws.Range("F370").Value = 3
$ws.Range("G370").Value = 0
$ws.Range("H370").Value = 0
$ws.Range("I370").Value = 10
$ws.Range("J370").Value = 45.5
$ws.Range("K370").Value = 44.5
$ws.Range("L370").Value = 190
$ws.Range("M370").Value = 0.5
$ws.Rows(371).ClearFormats()
$ws.Range("A371").Value = 1204
$ws.Range("B371").Value = 458
$ws.Range("C371").Value = "H"
$ws.Range("D371").Value = "Philadelphia"
$ws.Range("E371").Value = 7
$ws.Range("F371").Value = 14
$ws.Range("G371").Value = 7
$ws.Range("H371").Value = 7
$ws.Range("I371").Value = 35
$ws.Range("J371").Value = "7ev"
$ws.Range("K371").Value = 4.5
$ws.Range("L371").Value = -220
$ws.Range("M371").Value = 22.5
$ws.Rows(372).ClearFormats()
$ws.Range("A372").Value = 1204
$ws.Range("B372").Value = 459
$ws.Range("C372").Value = "V"
$ws.Range("D372").Value = "Jacksonville"
$ws.Range("E372").Value = 3
$ws.Range("F372").Value = 3
$ws.Range("G372").Value = 8
$ws.Range("H372").Value = 0
$ws.Range("I372").Value = 14
$ws.Range("J372").Value = 49.5
$ws.Range("K372").Value = 51
$ws.Range("L372").Value = 100
$ws.Range("M372").Value = 23.5
$ws.Rows(373).ClearFormats()
$ws.Range("A373").Value = 1204
$ws.Range("B373").Value = 460
$ws.Range("C373").Value = "H"
$ws.Range("D373").Value = "Detroit"
$ws.Range("E373").Value = 14
$ws.Range("F373").Value = 9
$ws.Range("G373").Value = 7
$ws.Range("H373").Value = 10
$ws.Range("I373").Value = 40
$ws.Range("J373").Value = 1
$ws.Range("K373").Value = "pk"
$ws.Range("L373").Value = -120
$ws.Range("M373").Value = 1.5
$ws.Rows(374).ClearFormats()
$ws.Range("A374").Value = 1204
$ws.Range("B374").Value = 461
$ws.Range("C374").Value = "V"
$ws.Range("D374").Value = "Washington"
$ws.Range("E374").Value = 10
$ws.Range("F374").Value = 3
$ws.Range("G374").Value = 0
$ws.Range("H374").Value = 7
$ws.Range("I374").Value = 20
$ws.Range("J374").Value = 1
$ws.Range("K374").Value = 1
$ws.Range("L374").Value = -135
$ws.Range("M374").Value = 2.5
$ws.Rows(375).ClearFormats()
$ws.Range("A375").Value = 1204
$ws.Range("B375").Value = 462
$ws.Range("C375").Value = "H"
$ws.Range("D375").Value = "NYGiants"
$ws.Range("E375").Value = 0
$ws.Range("F375").Value = 13
$ws.Range("G375").Value = 7
$ws.Range("H375").Value = 0
$ws.Range("I375").Value = 20
$ws.Range("J375").Value = 42
$ws.Range("K375").Value = 40.5
$ws.Range("L375").Value = 115
$ws.Range("M375").Value = 20.5
$ws.Rows(376).ClearFormats()
$ws.Range("A376").Value = 1204
$ws.Range("B376").Value = 463
$ws.Range("C376").Value = "V"
$ws.Range("D376").Value = "Cleveland"
$ws.Range("E376").Value = 0
$ws.Range("F376").Value = 7
$ws.Range("G376").Value = 7
$ws.Range("H376").Value = 13
$ws.Range("I376").Value = 27
$ws.Range("J376").Value = 5.5
$ws.Range("K376").Value = 9
$ws.Range("L376").Value = -350
$ws.Range("M376").Value = 4.5
$ws.Rows(377).ClearFormats()
$ws.Range("A377").Value = 1204
$ws.Range("B377").Value = 464
$ws.Range("C377").Value = "H"
$ws.Range("D377").Value = "Houston"
$ws.Range("E377").Value = 3
$ws.Range("F377").Value = 2
$ws.Range("G377").Value = 3
$ws.Range("H377").Value = 6
$ws.Range("I377").Value = 14
$ws.Range("J377").Value = 44
$ws.Range("K377").Value = 46
$ws.Range("L377").Value = 290
$ws.Range("M377").Value = 21
$ws.Rows(378).ClearFormats()
$ws.Range("A378").Value = 1204
$ws.Range("B378").Value = 465
$ws.Range("C378").Value = "V"
$ws.Range("D378").Value = "GreenBay"
$ws.Range("E378").Value = 0
$ws.Range("F378").Value = 10
$ws.Range("G378").Value = 0
$ws.Range("H378").Value = 18
$ws.Range("I378").Value = 28
$ws.Range("J378").Value = 2.5
$ws.Range("K378").Value = 3.5
$ws.Range("L378").Value = -185
$ws.Range("M378").Value = 3
$ws.Rows(379).ClearFormats()
$ws.Range("A379").Value = 1204
$ws.Range("B379").Value = 466
$ws.Range("C379").Value = "H"
$ws.Range("D379").Value = "Chicago"
$ws.Range("E379").Value = 10
$ws.Range("F379").Value = 6
$ws.Range("G379").Value = 3
$ws.Range("H379").Value = 0
$ws.Range("I379").Value = 19
$ws.Range("J379").Value = 44.5
$ws.Range("K379").Value = 44.5
$ws.Range("L379").Value = 165
$ws.Range("M379").Value = 23.5
$ws.Rows(380).ClearFormats()
$ws.Range("A380").Value = 1204
$ws.Range("B380").Value = 467
$ws.Range("C380").Value = "V"
$ws.Range("D380").Value = "Seattle"
$ws.Range("E380").Value = 7
$ws.Range("F380").Value = 7
$ws.Range("G380").Value = 3
$ws.Range("H380").Value = 10
$ws.Range("I380").Value = 27
$ws.Range("J380").Value = 3
$ws.Range("K380").Value = 7.5
$ws.Range("L380").Value = -290
$ws.Range("M380").Value = 3.5
$ws.Rows(381).ClearFormats()
$ws.Range("A381").Value = 1204
$ws.Range("B381").Value = 468
$ws.Range("C381").Value = "H"
$ws.Range("D381").Value = "LARams"
$ws.Range("E381").Value = 10
$ws.Range("F381").Value = 3
$ws.Range("G381").Value = 0
$ws.Range("H381").Value = 10
$ws.Range("I381").Value = 23
$ws.Range("J381").Value = 43
$ws.Range("K381").Value = 41
$ws.Range("L381").Value = 245
$ws.Range("M381").Value = 21
$ws.Rows(382).ClearFormats()
$ws.Range("A382").Value = 1204
$ws.Range("B382").Value = 469
$ws.Range("C382").Value = "V"
$ws.Range("D382").Value = "Miami"
$ws.Range("E382").Value = 7
$ws.Range("F382").Value = 3
$ws.Range("G382").Value = 0
$ws.Range("H382").Value = 7
$ws.Range("I382").Value = 17
$ws.Range("J382").Value = 46
$ws.Range("K382").Value = 46.5
$ws.Range("L382").Value = 190
$ws.Range("M382").Value = 0.5
$ws.Rows(383).ClearFormats()
$ws.Range("A383").Value = 1204
$ws.Range("B383").Value = 470
$ws.Range("C383").Value = "H"
$ws.Range("D383").Value = "SanFrancisco"
$ws.Range("E383").Value = 10
$ws.Range("F383").Value = 7
$ws.Range("G383").Value = 6
$ws.Range("H383").Value = 10
$ws.Range("I383").Value = 33
$ws.Range("J383").Value = 4.5
$ws.Range("K383").Value = 5
$ws.Range("L383").Value = -220
$ws.Range("M383").Value = 21
$ws.Rows(384).ClearFormats()
$ws.Range("A384").Value = 1204
$ws.Range("B384").Value = 471
$ws.Range("C384").Value = "V"
$ws.Range("D384").Value = "LAChargers"
$ws.Range("E384").Value = 7
$ws.Range("F384").Value = 6
$ws.Range("G384").Value = 0
$ws.Range("H384").Value = 7
$ws.Range("I384").Value = 20
$ws.Range("J384").Value = 3
$ws.Range("K384").Value = 49.5
$ws.Range("L384").Value = 120
$ws.Range("M384").Value = 24
$ws.Rows(385).ClearFormats()
$ws.Range("A385").Value = 1204
$ws.Range("B385").Value = 472
$ws.Range("C385").Value = "H"
$ws.Range("D385").Value = "LasVegas"
$ws.Range("E385").Value = 0
$ws.Range("F385").Value = 10
$ws.Range("G385").Value = 14
$ws.Range("H385").Value = 3
$ws.Range("I385").Value = 27
$ws.Range("J385").Value = 48.5
$ws.Range("K385").Value = 2.5
$ws.Range("L385").Value = -140
$ws.Range("M385").Value = 1
$ws.Rows(386).ClearFormats()
$ws.Range("A386").Value = 1204
$ws.Range("B386").Value = 473
$ws.Range("C386").Value = "V"
$ws.Range("D386").Value = "KansasCity"
$ws.Range("E386").Value = 3
$ws.Range("F386").Value = 7
$ws.Range("G386").Value = 14
$ws.Range("H386").Value = 0
$ws.Range("I386").Value = 24
$ws.Range("J386").Value = 2.5
$ws.Range("K386").Value = 2.5
$ws.Range("L386").Value = -140
$ws.Range("M386").Value = 2.5
$ws.Rows(387).ClearFormats()
$ws.Range("A387").Value = 1204
$ws.Range("B387").Value = 474
$ws.Range("C387").Value = "H"
$ws.Range("D387").Value = "Cincinnati"
$ws.Range("E387").Value = 7
$ws.Range("F387").Value = 7
$ws.Range("G387").Value = 3
$ws.Range("H387").Value = 10
$ws.Range("I387").Value = 27
$ws.Range("J387").Value = 51
$ws.Range("K387").Value = 53.5
$ws.Range("L387").Value = 120
$ws.Range("M387").Value = 27
$ws.Rows(388).ClearFormats()
$ws.Range("A388").Value = 1204
$ws.Range("B388").Value = 475
$ws.Range("C388").Value = "V"
$ws.Range("D388").Value = "Indianapolis"
$ws.Range("E388").Value = 10
$ws.Range("F388").Value = 3
$ws.Range("G388").Value = 6
$ws.Range("H388").Value = 0
$ws.Range("I388").Value = 19
$ws.Range("J388").Value = 45
$ws.Range("K388").Value = 44.5
$ws.Range("L388").Value = 400
$ws.Range("M388").Value = 21.5
$ws.Rows(389).ClearFormats()
$ws.Range("A389").Value = 1204
$ws.Range("B389").Value = 476
$ws.Range("C389").Value = "H"
$ws.Range("D389").Value = "Dallas"
$ws.Range("E389").Value = 7
$ws.Range("F389").Value = 14
$ws.Range("G389").Value = 0
$ws.Range("H389").Value = 33
$ws.Range("I389").Value = 54
$ws.Range("J389").Value = 9
$ws.Range("K389").Value = 11
$ws.Range("L389").Value = -500
$ws.Range("M389").Value = 3.5
